$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Burn-Down" estimate dropped by 10 hours starting 3/28 (row 28).
# C28 previously just carried the prior day's value forward ("=C27"); give
# it its own formula that knocks 10 hours off C27. Every later day in the
# column (C29:C42) still just carries the previous day's value forward, so
# they ripple down to the new, lower total automatically.
$ws.Range("C28").Formula = "=C27-10"

# Leave the selection where this edit was made, matching the saved view.
$ws.Range("C28").Select()
